# Lista de tareas pendientes API — add two new task rows (ratings + "Ponerlo
# Online") and mark the "Crear la segunda parte..." task's rating cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("Crear la segunda parte con la consulta de los datos de las pelis")
# now has its Completado flag set, and the description cell gets wrap text
# turned on (matches the wrapped layout used by the two new rows below it).
$ws.Range("D9").Value = 1
$ws.Range("B9").WrapText = $true

# Row 10: new task "Mejorar el diseño" / Priority "Alta"
$ws.Rows.Item(10).RowHeight = 49.8
$ws.Range("B10").Value = "Mejorar el diseño"
$ws.Range("C10").Value = "Alta"

# Row 11: new task "Ponerlo Online" / Priority "Media"
$ws.Range("B11").Value = "Ponerlo Online"
$ws.Range("C11").Value = "Media"

# Move the active selection to E11 (where the cursor ended up after the edits)
$ws.Range("E11").Select()
